$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "MISSING: $find"
    }
}

# --- KNN intro paragraph ---
Replace-Text "Effectively manage the daily maximum energy used, it is very important to analyze the energy uses in different sessions because the weather condition has a significant effect on" "Effectively managing the maximum daily energy used, it is very important to analyze the energy uses in different sessions because the weather significantly affects"

Replace-Text "utilized K-Nearest Neighbors (kNN) classifier for daily maximum energy uses prediction" "utilized K-Nearest Neighbors (kNN) classifier for maximum daily energy uses prediction"

Replace-Text "K-Nearest Neighbors is one of the simplest supervised classification techniques that provided good performance results for" "K-Nearest Neighbors is one of the most straightforward supervised classification techniques that provide exemplary performance results for"

# --- Model evaluation paragraph ---
Replace-Text "we have considered four stages namely data processing" "we have considered four stages: data processing"

Replace-Text "The providing datasets containing weather conditions" "The datasets containing weather conditions"

Replace-Text "There are two steps are involved in the prediction stage namely the training stage and the testing stage. In the training stage," "Two steps are involved in the prediction stage, namely the training and testing stages. In the training stage,"

Replace-Text "During the testing phase, the KNN classifier is given unlabelled data points and the algorithm generates" "During testing phase, the KNN classifier is given unlabelled data points, and the algorithm generates"

# --- First model accuracy paragraph ---
Replace-Text "We have split the train and test size into 80% and 20% in 42 random states. The highest accuracy has been observed for our prediction is 72.7%" "We have split the train and test size into 80% and 20% in 42 random states. The highest accuracy observed for our prediction is 72.7%"

# --- Second model accuracy paragraph ---
Replace-Text "We have split the train and test size into 87% and 13% in 42 random states. The highest accuracy has been observed for our prediction is 55.1%" "We have split the train and test size into 87% and 13% in 42 random states. The highest accuracy observed for our prediction is 55.1%"

# --- K-Fold body paragraph ---
Replace-Text "K number of folds, in which each fold is utilized as a testing set in a particular selective point and the remaining folds" "K number of folds. Each fold is utilized as testing set in a particular selective point, and the remaining folds"

Replace-Text "the accuracy of the model performance. In both models, we have applied the K-Fold method to understand the model performance in a particular training set. First of all, we have split the data set into 10 folds and shuffled each of the folds in every test in 42 random states." "the accuracy of the model performance. We have applied the K-Fold method in both models to understand the model performance in a particular training set. First, we split the data set into ten folds and shuffled each of the folds in every test in 42 random states."

# --- First K-Fold model accuracy paragraph ---
Replace-Text "For the first model, we have trained the feature dataset. The model accuracy has been observed for our prediction is 49.9%" "For the first model, we have trained the feature dataset. The model accuracy observed for our prediction is 49.9%"

# --- Second K-Fold model accuracy paragraph ---
Replace-Text "For the second model, we have trained the feature dataset. The model accuracy has been observed for our prediction is 49.9%" "For the second model, we have trained the feature dataset. The model accuracy observed for our prediction is 49.9%"

Write-Output "done"
